$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: replace the "E" difficulty marker with a numeric count
$ws.Range("D2").Value = 10

# Rows 4-14: fill in the "Total / Day" numeric counts in column D
$ws.Range("D4").Value = 5
$ws.Range("D5").Value = 4
$ws.Range("D6").Value = 3
$ws.Range("D7").Value = 13
$ws.Range("D8").Value = 5
$ws.Range("D9").Value = 3
$ws.Range("D10").Value = 4
$ws.Range("D11").Value = 5
$ws.Range("D12").Value = 7
$ws.Range("D13").Value = 10
$ws.Range("D14").Value = 5

# Row 1: replace the "Difficulty" header with the running total formula
$ws.Range("D1").Formula = "=SUM(D2:D15)"

# New log entries (Jogging / Check / new date ranges / new notes)
$ws.Range("E15").Value = "Jogging"
$ws.Range("H9").Value = "July 19 - July 21 "

# G9 keeps its original "quote-prefixed" cell format even though the
# text changes, so stash the format on a scratch cell, overwrite the
# value, then restore the format.
$ws.Range("G9").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null
$ws.Range("G9").Value = "Check"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("G9").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Clear() | Out-Null

# H10 picks up H9's date-range cell format
$ws.Range("H9").Copy() | Out-Null
$ws.Range("H10").PasteSpecial(-4122) | Out-Null
$ws.Range("H10").Value = "July 21 - July  "

$ws.Range("C15").Value = "P118, P347, P380, P498, P724, P747`nHashTable finished, Array and String Started"
$ws.Range("C15").WrapText = $true

# Row 15: finish populating the new row (previously blank)
$ws.Range("A15").Value = 44033
$ws.Range("D15").Value = 6
$ws.Rows.Item(15).RowHeight = 30

# Update the active selection to match the new working cell
[void]$ws.Range("C15").Select()
